$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 552.6667
$ws.Range("I2").Value = 400
$ws.Range("J2").Value = 583.2
$ws.Range("K2").Value = 400
$ws.Range("L2").Value = 583.2
$ws.Range("M2").Value = -287
$ws.Range("N2").Value = -809.2
$ws.Range("H12").Value = 140.4
$ws.Range("I12").Value = 200
$ws.Range("J12").Value = 51
$ws.Range("K12").Value = 200
$ws.Range("L12").Value = 51
$ws.Range("M12").Value = -30
$ws.Range("N12").Value = -391
$ws.Range("H69").Value = 4040.9375
$ws.Range("I69").Value = 3593.3333
$ws.Range("J69").Value = 4144.231
$ws.Range("K69").Value = 10779.9999
$ws.Range("L69").Value = 12432.693
$ws.Range("M69").Value = -9905.999899999999
$ws.Range("N69").Value = -14180.693
$ws.Range("H72").Value = 4040.9375
$ws.Range("I72").Value = 3593.3333
$ws.Range("J72").Value = 4144.231
$ws.Range("K72").Value = 32339.9997
$ws.Range("L72").Value = 37298.079
$ws.Range("M72").Value = -27971.9997
$ws.Range("N72").Value = -46034.079
$ws.Range("H86").Value = 8610.412
$ws.Range("I86").Value = 12714.9
$ws.Range("J86").Value = 2746.8572
$ws.Range("K86").Value = 12714.9
$ws.Range("L86").Value = 2746.8572
$ws.Range("M86").Value = -11591.9
$ws.Range("N86").Value = -4992.8572
$ws.Range("H89").Value = 8610.412
$ws.Range("I89").Value = 12714.9
$ws.Range("J89").Value = 2746.8572
$ws.Range("K89").Value = 63574.5
$ws.Range("L89").Value = 13734.286
$ws.Range("M89").Value = -57958.5
$ws.Range("N89").Value = -24966.286
$ws.Range("H107").Value = 317.14285
$ws.Range("I107").Value = 213.6
$ws.Range("J107").Value = 576
$ws.Range("K107").Value = 213.6
$ws.Range("L107").Value = 576
$ws.Range("M107").Value = 1706.4
$ws.Range("N107").Value = -4416
$ws.Range("H112").Value = 1158313
$ws.Range("J112").Value = 1544292.4
$ws.Range("L112").Value = 4632877.199999999
$ws.Range("N112").Value = -4635093.199999999
$ws.Range("H113").Value = 2891.8076
$ws.Range("I113").Value = 2517
$ws.Range("J113").Value = 3166.6667
$ws.Range("K113").Value = 2517
$ws.Range("L113").Value = 3166.6667
$ws.Range("M113").Value = 737
$ws.Range("N113").Value = -9674.6667
$ws.Range("H127").Value = 899.875
$ws.Range("I127").Value = 633.6667
$ws.Range("J127").Value = 1698.5
$ws.Range("K127").Value = 1901.0001
$ws.Range("L127").Value = 5095.5
$ws.Range("M127").Value = 3058.9999
$ws.Range("N127").Value = -15015.5
$ws.Range("H129").Value = 7784.25
$ws.Range("I129").Value = 446
$ws.Range("J129").Value = 9838.959999999999
$ws.Range("K129").Value = 1338
$ws.Range("L129").Value = 29516.88
$ws.Range("M129").Value = 3662
$ws.Range("N129").Value = -39516.88

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 3174.3125
$ws.Range("I97").Value = 4639.778
$ws.Range("K97").Value = 4639.778
$ws.Range("M97").Value = -4143.778

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4102.8125
$ws.Range("I20").Value = 1667
$ws.Range("J20").Value = 9461.6
$ws.Range("K20").Value = 1667
$ws.Range("L20").Value = 9461.6
$ws.Range("M20").Value = -1420
$ws.Range("N20").Value = -9955.6
$ws.Range("H82").Value = 14107.533
$ws.Range("I82").Value = 4591.3
$ws.Range("J82").Value = 33140
$ws.Range("K82").Value = 4591.3
$ws.Range("L82").Value = 33140
$ws.Range("M82").Value = -4208.3
$ws.Range("N82").Value = -33906
$ws.Range("H85").Value = 14107.533
$ws.Range("I85").Value = 4591.3
$ws.Range("J85").Value = 33140
$ws.Range("K85").Value = 4591.3
$ws.Range("L85").Value = 33140
$ws.Range("M85").Value = -3265.3
$ws.Range("N85").Value = -35792

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 4199.4375
$ws.Range("I99").Value = 4725.4546
$ws.Range("J99").Value = 3042.2
$ws.Range("K99").Value = 4725.4546
$ws.Range("L99").Value = 3042.2
$ws.Range("M99").Value = -3227.4546
$ws.Range("N99").Value = -6038.2
$ws.Range("H126").Value = 4199.4375
$ws.Range("I126").Value = 4725.4546
$ws.Range("J126").Value = 3042.2
$ws.Range("K126").Value = 14176.3638
$ws.Range("L126").Value = 9126.599999999999
$ws.Range("M126").Value = -11706.3638
$ws.Range("N126").Value = -14066.6
$ws.Range("H134").Value = 2962.4546
$ws.Range("I134").Value = 3368
$ws.Range("J134").Value = 1137.5
$ws.Range("K134").Value = 10104
$ws.Range("L134").Value = 3412.5
$ws.Range("M134").Value = -7569
$ws.Range("N134").Value = -8482.5
$ws.Range("H135").Value = 9335086
$ws.Range("J135").Value = 9335086
$ws.Range("L135").Value = 9335086
$ws.Range("N135").Value = -9345226

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1040.0605
$ws.Range("I5").Value = 926.3125
$ws.Range("J5").Value = 1147.1177
$ws.Range("K5").Value = 2778.9375
$ws.Range("L5").Value = 3441.3531
$ws.Range("M5").Value = -2666.9375
$ws.Range("N5").Value = -3665.3531
$ws.Range("H63").Value = 7700.857
$ws.Range("I63").Value = 1203
$ws.Range("J63").Value = 10300
$ws.Range("K63").Value = 3609
$ws.Range("L63").Value = 30900
$ws.Range("M63").Value = -2860
$ws.Range("N63").Value = -32398
$ws.Range("H66").Value = 7700.857
$ws.Range("I66").Value = 1203
$ws.Range("J66").Value = 10300
$ws.Range("K66").Value = 10827
$ws.Range("L66").Value = 92700
$ws.Range("M66").Value = -7083
$ws.Range("N66").Value = -100188
$ws.Range("H82").Value = 4025.25
$ws.Range("I82").Value = 1750.75
$ws.Range("J82").Value = 6299.75
$ws.Range("K82").Value = 5252.25
$ws.Range("L82").Value = 18899.25
$ws.Range("M82").Value = -4846.25
$ws.Range("N82").Value = -19711.25
$ws.Range("H85").Value = 4025.25
$ws.Range("I85").Value = 1750.75
$ws.Range("J85").Value = 6299.75
$ws.Range("K85").Value = 5252.25
$ws.Range("L85").Value = 18899.25
$ws.Range("M85").Value = -3848.25
$ws.Range("N85").Value = -21707.25
$ws.Range("H100").Value = 12212.375
$ws.Range("J100").Value = 13800
$ws.Range("L100").Value = 41400
$ws.Range("N100").Value = -43022
$ws.Range("H109").Value = 3087.8
$ws.Range("I109").Value = 2446.6667
$ws.Range("J109").Value = 4049.5
$ws.Range("K109").Value = 7340.000100000001
$ws.Range("L109").Value = 12148.5
$ws.Range("M109").Value = -6300.000100000001
$ws.Range("N109").Value = -14228.5
$ws.Range("H115").Value = 4853.6895
$ws.Range("J115").Value = 3433.8076
$ws.Range("L115").Value = 10301.4228
$ws.Range("N115").Value = -12651.4228
$ws.Range("H118").Value = 8245.772000000001
$ws.Range("I118").Value = 1136.3334
$ws.Range("J118").Value = 9368.315000000001
$ws.Range("K118").Value = 3409.0002
$ws.Range("L118").Value = 28104.945
$ws.Range("M118").Value = -2166.0002
$ws.Range("N118").Value = -30590.945
$ws.Range("H131").Value = 907.5
$ws.Range("I131").Value = 512.8
$ws.Range("J131").Value = 1236.4166
$ws.Range("K131").Value = 1538.4
$ws.Range("L131").Value = 3709.2498
$ws.Range("M131").Value = 3501.6
$ws.Range("N131").Value = -13789.2498
$ws.Range("H132").Value = 886.44446
$ws.Range("I132").Value = 513
$ws.Range("J132").Value = 1633.3334
$ws.Range("K132").Value = 4617
$ws.Range("L132").Value = 14700.0006
$ws.Range("M132").Value = -2087
$ws.Range("N132").Value = -19760.0006
$ws.Range("H135").Value = 1040.0605
$ws.Range("I135").Value = 926.3125
$ws.Range("J135").Value = 1147.1177
$ws.Range("K135").Value = 8336.8125
$ws.Range("L135").Value = 10324.0593
$ws.Range("M135").Value = -5801.8125
$ws.Range("N135").Value = -15394.0593

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2514.2856
$ws.Range("I80").Value = 2575
$ws.Range("J80").Value = 2433.3333
$ws.Range("K80").Value = 2575
$ws.Range("L80").Value = 2433.3333
$ws.Range("M80").Value = -1577
$ws.Range("N80").Value = -4429.3333
$ws.Range("H83").Value = 2514.2856
$ws.Range("I83").Value = 2575
$ws.Range("J83").Value = 2433.3333
$ws.Range("K83").Value = 12875
$ws.Range("L83").Value = 12166.6665
$ws.Range("M83").Value = -7883
$ws.Range("N83").Value = -22150.6665

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2085.4375
$ws.Range("I40").Value = 1647.25
$ws.Range("J40").Value = 3400
$ws.Range("K40").Value = 1647.25
$ws.Range("L40").Value = 3400
$ws.Range("M40").Value = -1511.25
$ws.Range("N40").Value = -3672
$ws.Range("H92").Value = 28000
$ws.Range("J92").Value = 28000
$ws.Range("L92").Value = 28000
$ws.Range("N92").Value = -32992

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5478.5713
$ws.Range("I62").Value = 5666.6665
$ws.Range("J62").Value = 5337.5
$ws.Range("K62").Value = 5666.6665
$ws.Range("L62").Value = 5337.5
$ws.Range("M62").Value = -5042.6665
$ws.Range("N62").Value = -6585.5
$ws.Range("H65").Value = 5478.5713
$ws.Range("I65").Value = 5666.6665
$ws.Range("J65").Value = 5337.5
$ws.Range("K65").Value = 28333.3325
$ws.Range("L65").Value = 26687.5
$ws.Range("M65").Value = -25213.3325
$ws.Range("N65").Value = -32927.5
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()  # cell removed in target

